$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the target range to Text format so the new numeric-looking / percent-looking
# strings are stored as text (matching the original inlineStr "t" cells) rather than
# being auto-converted to numbers by Excel.
$fmtRange = $ws.Range("D2:E47")
$fmtRange.NumberFormat = "@"

$ws.Range("D2").Value = "275.96"
$ws.Range("E2").Value = "0.57%"
$ws.Range("E3").Value = "1.67%"
$ws.Range("D4").Value = "4.845"
$ws.Range("E4").Value = "-0.32%"
$ws.Range("E5").Value = "1.15%"
$ws.Range("D6").Value = "6.945"
$ws.Range("E6").Value = "0.79%"
$ws.Range("D7").Value = "1.205"
$ws.Range("E7").Value = "-5.69%"
$ws.Range("D8").Value = "0.8780"
$ws.Range("E8").Value = "0.69%"
$ws.Range("D9").Value = "0.1520"
$ws.Range("D10").Value = "0.05108"
$ws.Range("E10").Value = "1.24%"
$ws.Range("D11").Value = "0.07524"
$ws.Range("E11").Value = "2.21%"
$ws.Range("E12").Value = "1.70%"
$ws.Range("E13").Value = "-0.62%"
$ws.Range("D14").Value = "0.001571"
$ws.Range("E14").Value = "0.19%"
$ws.Range("D15").Value = "0.0006403"
$ws.Range("E15").Value = "1.52%"
$ws.Range("D16").Value = "0.006151"
$ws.Range("E16").Value = "2.25%"
$ws.Range("D17").Value = "3.482"
$ws.Range("E17").Value = "0.97%"
$ws.Range("D18").Value = "3.305"
$ws.Range("E18").Value = "-0.52%"
$ws.Range("E19").Value = "-1.92%"
$ws.Range("E20").Value = "0.29%"
$ws.Range("E21").Value = "1.85%"
$ws.Range("D22").Value = "3.914"
$ws.Range("E22").Value = "0.39%"
$ws.Range("D23").Value = "0.04402"
$ws.Range("E23").Value = "1.09%"
$ws.Range("D25").Value = "0.001178"
$ws.Range("E25").Value = "0.03%"
$ws.Range("D26").Value = "0.003860"
$ws.Range("E26").Value = "-9.51%"
$ws.Range("E27").Value = "0.11%"
$ws.Range("E28").Value = "14.75%"
$ws.Range("D40").Value = "0.04159"
$ws.Range("E40").Value = "2.95%"
$ws.Range("D41").Value = "0.006816"
$ws.Range("E41").Value = "2.84%"
$ws.Range("E42").Value = "0.51%"
$ws.Range("E43").Value = "2.97%"
$ws.Range("D44").Value = "0.01187"
$ws.Range("E44").Value = "-2.97%"
$ws.Range("D45").Value = "0.00005278"
$ws.Range("E45").Value = "-0.48%"
$ws.Range("D46").Value = "1.681"
$ws.Range("E46").Value = "16.02%"
$ws.Range("D47").Value = "0.01850"
$ws.Range("E47").Value = "-7.39%"

# Reset the style back to Normal so no stray NumberFormat/style metadata is left
# behind on the range (the cells keep their text values, just drop the applied "@" xf).
$fmtRange.Style = "Normal"
